{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Change described by the diff: the command-line paragraph that read\n//   \"python batchUpload.py\"\n// becomes\n//   \"python.exe batchUpload.py batchUploadGuest.json.gz\"\n// (a guest batch file is now referenced alongside the original order file).\n// The remaining diff hunks are just Word re-flowing/merging runs that used\n// to be split apart by <w:proofErr/> spell-check markers - the visible text\n// for those paragraphs is unchanged, so no action is required for them.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraph that contains the command line we need to update.\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"python batchUpload.py\") !== -1) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (target) {\n  // Work right-to-left so earlier ranges stay valid while later text grows.\n\n  // 1) \"batchUpload.py\" -> \"batchUpload.py batchUploadGuest.json.gz\"\n  const tailResults = target.search(\"batchUpload.py\", { matchCase: true });\n  tailResults.load(\"items\");\n  await context.sync();\n  if (tailResults.items.length > 0) {\n    tailResults.items[0].insertText(\n      \"batchUpload.py batchUploadGuest.json.gz\",\n      Word.InsertLocation.replace\n    );\n    await context.sync();\n  }\n\n  // 2) \"python\" -> \"python.exe\" (insert \".exe\" right after \"python\")\n  const pyResults = target.search(\"python\", { matchCase: true });\n  pyResults.load(\"items\");\n  await context.sync();\n  if (pyResults.items.length > 0) {\n    pyResults.items[0].insertText(\".exe\", Word.InsertLocation.after);\n    await context.sync();\n  }\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document.\n#\n# Change described by the diff: the command-line paragraph that read\n#   \"python batchUpload.py\"\n# becomes\n#   \"python.exe batchUpload.py batchUploadGuest.json.gz\"\n# (a guest batch file is now referenced alongside the original order file).\n# The remaining diff hunks are just Word re-flowing/merging runs that used\n# to be split apart by proofErr spell-check markers - the visible text for\n# those paragraphs is unchanged, so no action is required for them.\n\n$d = $word.ActiveDocument\n\n# Find the paragraph holding the command line we need to update.\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*python batchUpload.py*\") {\n        $target = $p.Range\n        break\n    }\n}\n\nif ($target -ne $null) {\n    $rng = $target.Duplicate\n    $found = $rng.Find.Execute(\"python batchUpload.py\", $true)\n    if ($found) {\n        $rng.Text = \"python.exe batchUpload.py batchUploadGuest.json.gz\"\n    }\n}\n"}
